# Split QO subset to ML + MB
#
# The "Easy Selector Jump Start" sheet listed a combined "Indoor 1Ph3W Main
# Lugs and Main Breaker" product (with an Attribute of either "Circuit
# Breaker" or "Lugs") across rows 2-86. This change splits that single
# product into two distinct catalog entries - "Indoor 1Ph3W Main Breaker"
# (rows 2-41) and "Indoor 1Ph3W Main Lugs" (rows 42-86) - each now carrying
# an Attribute of "n/a" since the breaker/lugs distinction is now baked
# into the product name itself.
#
# It also grows the _FilterDatabase defined name / AutoFilter range and
# resets the frozen-pane view back to the top of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Easy Selector Jump Start")

# Rows 2-41: "...Main Lugs and Main Breaker" / "Circuit Breaker"
#            -> "Indoor 1Ph3W Main Breaker" / "n/a"
$ws.Range("D2:D41").Value = "Indoor 1Ph3W Main Breaker"
$ws.Range("E2:E41").Value = "n/a"

# Rows 42-86: "...Main Lugs and Main Breaker" / "Lugs"
#             -> "Indoor 1Ph3W Main Lugs" / "n/a"
$ws.Range("D42:D86").Value = "Indoor 1Ph3W Main Lugs"
$ws.Range("E42:E86").Value = "n/a"

# Grow the hidden _xlnm._FilterDatabase defined name to cover the full table
$fdb = $wb.Names.Item("Easy Selector Jump Start!_FilterDatabase")
$fdb.RefersTo = "='Easy Selector Jump Start'!`$A`$1:`$H`$278"

# Re-apply the AutoFilter over the full data range A1:H278
$ws.AutoFilterMode = $False
$ws.Range("A1:H278").AutoFilter()

# Reset the frozen-pane view: scroll back to the top of the sheet and move
# the active-cell selection from the old bottom (D279) to D43
$ws.Activate()
$ws.Range("D43").Select()
$wb.Worksheets.Item("Configuration").Activate()
